# GoalConfig.xlsx — update Goal1 X/Y position values and move the
# selection, matching the authored diff:
#   B2: 132  -> 154   (msg.Pose.Position.X)
#   B3: -314 -> -339  (msg.Pose.Position.Y)
#   active cell/selection moves from B3 to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 154
$ws.Range("B3").Value = -339

# Move the active cell / selection to B4 (was B3).
$ws.Range("B4").Select()
